# Daily attendance processing - 2026-01-24 00:00:49
# Swap the order of the two comma-separated entries in the "Recorded By"
# column (G) for every session row, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Rows whose "Recorded By" list has a different shape (single entry, three
# entries, or includes "backup@backdoor.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -and $value -is [string] -and $value.Contains(", ")) {
        $parts = $value -split ", "

        if ($parts.Count -eq 2 -and -not $value.Contains("backup@backdoor.com")) {
            $cell.Value2 = $parts[1] + ", " + $parts[0]
        }
    }
}
